$d = $word.ActiveDocument

function Replace-ParagraphText($paragraphIndex, $oldText, $newText) {
    $p = $d.Paragraphs($paragraphIndex)
    $rng = $p.Range
    $find = $rng.Find
    $find.ClearFormatting()
    $find.Replacement.ClearFormatting()
    $result = $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
    if (-not $result) {
        Write-Host "WARNING: replace failed for paragraph $paragraphIndex : $oldText"
    }
}

# ---------------------------------------------------------------------------
# 1. Objetivos paragraph: old Objetivos text -> old "Programa resumido" text
# ---------------------------------------------------------------------------
Replace-ParagraphText 6 `
    "Apresentar os princípios básicos da Seleção de Materiais para aplicação em Engenharia." `
    "Aspectos gerais e critérios de seleção de materiais estruturais. Aspectos dos principais mecanismos de falha em componentes estruturais. Seleção de materiais e análise para diferentes modos de carregamento. Seleção de materiais sob diferentes condições de temperatura. Materiais resistentes à corrosão e oxidação. Tribologia: atrito e desgaste. Tratamentos superficiais."

# ---------------------------------------------------------------------------
# 2. Docente paragraph: old Docente text -> old Objetivos text
# ---------------------------------------------------------------------------
Replace-ParagraphText 8 `
    "5840622 - Miguel Justino Ribeiro Barboza" `
    "Apresentar os princípios básicos da Seleção de Materiais para aplicação em Engenharia."

# ---------------------------------------------------------------------------
# 3. "Programa resumido" paragraph: old text -> old "Programa" text
# ---------------------------------------------------------------------------
Replace-ParagraphText 10 `
    "Aspectos gerais e critérios de seleção de materiais estruturais. Aspectos dos principais mecanismos de falha em componentes estruturais. Seleção de materiais e análise para diferentes modos de carregamento. Seleção de materiais sob diferentes condições de temperatura. Materiais resistentes à corrosão e oxidação. Tribologia: atrito e desgaste. Tratamentos superficiais." `
    "1. Principais mecanismos de falha em componentes estruturais: efeitos do meio e temperatura. Critérios de falha. 2. Seleção de materiais para aplicações sob a ação de cargas estáticas. Materiais metálicos, cerâmicos, poliméricos e compósitos. 3. Seleção de materiais para aplicações sob a ação de cargas dinâmicas: O fenômeno da fadiga e efeitos da presença de entalhes em componentes mecânicos. 4. Critérios de seleção de materiais para aplicações em temperaturas elevadas. O fenômeno da fluência e a tolerância ao dano. Seleção de materiais para alta temperatura. Aços especiais, superligas, materiais cerâmicos e compósitos. 5. Materiais para temperaturas criogênicas. A transição dúctil-frágil. 6. Aspectos fundamentais do estudo de tribologia: desgaste, atrito e tratamentos superficiais. 7. Fundamentos, seleção e proteção contra oxidação. 8. Seleção de materiais em meios corrosivos. Corrosão sob tensão."

# ---------------------------------------------------------------------------
# 4. "Programa" paragraph: old text -> old "Método" intro text
# ---------------------------------------------------------------------------
Replace-ParagraphText 12 `
    "1. Principais mecanismos de falha em componentes estruturais: efeitos do meio e temperatura. Critérios de falha. 2. Seleção de materiais para aplicações sob a ação de cargas estáticas. Materiais metálicos, cerâmicos, poliméricos e compósitos. 3. Seleção de materiais para aplicações sob a ação de cargas dinâmicas: O fenômeno da fadiga e efeitos da presença de entalhes em componentes mecânicos. 4. Critérios de seleção de materiais para aplicações em temperaturas elevadas. O fenômeno da fluência e a tolerância ao dano. Seleção de materiais para alta temperatura. Aços especiais, superligas, materiais cerâmicos e compósitos. 5. Materiais para temperaturas criogênicas. A transição dúctil-frágil. 6. Aspectos fundamentais do estudo de tribologia: desgaste, atrito e tratamentos superficiais. 7. Fundamentos, seleção e proteção contra oxidação. 8. Seleção de materiais em meios corrosivos. Corrosão sob tensão." `
    "Este curso deverá conter duas avaliações escritas denominadas P1 e P2. A P2 deverá englobar toda a matéria ministrada ao longo do semestre, abrangendo todos os tópicos previstos na ementa."

# ---------------------------------------------------------------------------
# 5. "Avaliação" paragraph (index 14): restructure runs.
#    a) Remove the "Este curso deverá ... ementa." run (+ its trailing break)
#       together with the bold "Critério: " label that followed it.
#    b) Rename the bold "Norma de recuperação: " label (the first one, which
#       precedes the "A recuperação..." text) to "Critério: ".
#    c) Append a line break, a new bold "Norma de recuperação: " label and
#       the old Bibliografia text at the end of the paragraph.
# ---------------------------------------------------------------------------
$vtab = [char]11

# 5a. remove the obsolete "Método" explanation + following "Critério: " label
$p14 = $d.Paragraphs(14)
$find5a = $p14.Range.Find
$find5a.ClearFormatting()
$find5a.Replacement.ClearFormatting()
$removeText = "Este curso deverá conter duas avaliações escritas denominadas P1 e P2. A P2 deverá englobar toda a matéria ministrada ao longo do semestre, abrangendo todos os tópicos previstos na ementa." + $vtab + "Critério: "
$find5a.Execute($removeText, $true, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null

# 5b. rename the remaining bold "Norma de recuperação: " label to "Critério: "
$p14b = $d.Paragraphs(14)
$find5b = $p14b.Range.Find
$find5b.ClearFormatting()
$find5b.Replacement.ClearFormatting()
$find5b.Replacement.Font.Bold = $true
$find5b.Execute("Norma de recuperação: ", $true, $false, $false, $false, $false, $true, 1, $false, "Critério: ", 2) | Out-Null

# 5c. add trailing break after "MF=(M+RC)/2"
$p14c = $d.Paragraphs(14)
$find5c = $p14c.Range.Find
$find5c.ClearFormatting()
$find5c.Replacement.ClearFormatting()
$find5c.Execute("MF=(M+RC)/2", $true, $false, $false, $false, $false, $true, 1, $false, "MF=(M+RC)/2" + $vtab, 2) | Out-Null

# 5d. append new bold "Norma de recuperação: " label at the very end of the paragraph
$endPos = $d.Paragraphs(14).Range.End
$insertLabel = $d.Range($endPos - 1, $endPos - 1)
$insertLabel.InsertAfter("Norma de recuperação: ")
$labelStart = $endPos - 1
$labelEnd = $d.Paragraphs(14).Range.End - 1
$boldRng = $d.Range($labelStart, $labelEnd)
$boldRng.Font.Bold = $true

# 5e. append old Bibliografia text (plain) right after the new label
$biblioText = "1.Ashby, M. F. Materials Selection in Mechanical Design, Butterworth, Oxford, 2005. 2. ASM Metals Handbook - Properties and Selection: Irons, Steels and High - Performance Alloys - v.1 - 1990. 3. ASM Metals Handbook - Properties and Selection: Nonferrous Alloys and Special - Purpose Materials - v.2 - 1990. 4. Meyers, M.; Chawla, K. Mechanical Behavior of Materials. Ed. Cambridge University Press, 2009. 5. Van Vlack, L.H., Propriedades dos Materiais Cerâmicos. Ed. Edgard Blücher Ltda., 1973. 6. Dowling, E. M. Mechanical behavior of materials: engineering methods for deformation, fracture and fatigue. New Jersey, Prentice Hall, 1999. 7. Biasotto, E., Polímeros como Materiais de Engenharia. Ed. Edgard Blücher Ltda., 1991. 8. Rosen, S.L., Fundamental Principles of Polymeric Materials. Ed. John Wiley & Sons, Inc., 1993. 9. Bhushan, B. Introduction to Tribology, 2nd Edition, John Wiley & Sons. 2013. 10. Roberge, P. R. Corrosion engineering: principles and practice. The McGraw-Hill Companies, Inc., 2008. 11. Gentil, V. Corrosão, Ed. LTC, 2011. 12. Crane, F.A., Charles, J.A., Selection of Engineering Materials, Butterworth, 1984. 13. Chiaverini, V., Aços e Ferros Fundidos, Associação Brasileira de Materiais - ABM, São Paulo, 1988. 14. Reed, R. C. The superalloys: fundamentals and applications. Ed. Cambridge, USA, 2006."
$endPos2 = $d.Paragraphs(14).Range.End
$insertBiblio = $d.Range($endPos2 - 1, $endPos2 - 1)
$insertBiblio.InsertAfter($biblioText)

# ---------------------------------------------------------------------------
# 6. Bibliografia paragraph: old bibliography text -> old Docente text
# ---------------------------------------------------------------------------
Replace-ParagraphText 16 `
    "1.Ashby, M. F. Materials Selection in Mechanical Design, Butterworth, Oxford, 2005. 2. ASM Metals Handbook - Properties and Selection: Irons, Steels and High - Performance Alloys - v.1 - 1990. 3. ASM Metals Handbook - Properties and Selection: Nonferrous Alloys and Special - Purpose Materials - v.2 - 1990. 4. Meyers, M.; Chawla, K. Mechanical Behavior of Materials. Ed. Cambridge University Press, 2009. 5. Van Vlack, L.H., Propriedades dos Materiais Cerâmicos. Ed. Edgard Blücher Ltda., 1973. 6. Dowling, E. M. Mechanical behavior of materials: engineering methods for deformation, fracture and fatigue. New Jersey, Prentice Hall, 1999. 7. Biasotto, E., Polímeros como Materiais de Engenharia. Ed. Edgard Blücher Ltda., 1991. 8. Rosen, S.L., Fundamental Principles of Polymeric Materials. Ed. John Wiley & Sons, Inc., 1993. 9. Bhushan, B. Introduction to Tribology, 2nd Edition, John Wiley & Sons. 2013. 10. Roberge, P. R. Corrosion engineering: principles and practice. The McGraw-Hill Companies, Inc., 2008. 11. Gentil, V. Corrosão, Ed. LTC, 2011. 12. Crane, F.A., Charles, J.A., Selection of Engineering Materials, Butterworth, 1984. 13. Chiaverini, V., Aços e Ferros Fundidos, Associação Brasileira de Materiais - ABM, São Paulo, 1988. 14. Reed, R. C. The superalloys: fundamentals and applications. Ed. Cambridge, USA, 2006." `
    "5840622 - Miguel Justino Ribeiro Barboza"

Write-Host "All edits applied."
